$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number format codes used by the existing color-coded columns (D:H)
$fmt164 = "_(* #,##0.0000_);_(* \(#,##0.0000\);_(* ""-""????_);_(@_)"
$fmt41  = "_(* #,##0_);_(* \(#,##0\);_(* ""-""_);_(@_)"

# Existing theme colors used for columns D-H (Red, Orange, Yellow, Green, Blue)
$colRed    = 0x0000FF
$colOrange = 0x00C0FF
$colYellow = 0x00FFFF
$colGreen  = 0x50D092
$colBlue   = 0xF0B000

# --- Row 16 (written first so "weighted sum" / "binary-classifier score" become
#     shared-string indices 24 / 25, matching the reference order) ---
$ws.Range("C16").Value = "weighted sum"
$ws.Range("B16").Value = "binary-classifier score"

# New cell styles (numFmtId 41 + fill) are created in the same order as the
# reference file so the resulting style indices (10-14) line up exactly.
$ws.Range("H16").NumberFormat = $fmt41
$ws.Range("H16").Interior.Color = $colRed
$ws.Range("H16").Formula = "=2*5+3*4+0*3+2*2+0*1"

$ws.Range("D16").NumberFormat = $fmt41
$ws.Range("D16").Interior.Color = $colOrange
$ws.Range("D16").Formula = "=3*5+2*4+0*3+0*2+2*1"

$ws.Range("G16").NumberFormat = $fmt41
$ws.Range("G16").Interior.Color = $colYellow
$ws.Range("G16").Formula = "=0*5+0*4+7*3+0*2+0*1"

$ws.Range("E16").NumberFormat = $fmt41
$ws.Range("E16").Interior.Color = $colGreen
$ws.Range("E16").Formula = "=0*5+2*4+0*3+5*2+0*1"

$ws.Range("F16").NumberFormat = $fmt41
$ws.Range("F16").Interior.Color = $colBlue
$ws.Range("F16").Formula = "=2*5+0*4+0*3+0*2+5*1"

# --- Row 14 ---
$ws.Range("B14").Value = "lower bound of CI for optimal F1 measure"
$ws.Range("C14").Value = "optimal F1 measure - critical t value * standard error of optimal F1 measure"

$ws.Range("D14").NumberFormat = $fmt164
$ws.Range("D14").Interior.Color = $colOrange
$ws.Range("D14").Value = 0.89500000000000002

$ws.Range("E14").NumberFormat = $fmt164
$ws.Range("E14").Interior.Color = $colGreen
$ws.Range("E14").Value = 0.8841

$ws.Range("F14").NumberFormat = $fmt164
$ws.Range("F14").Interior.Color = $colBlue
$ws.Range("F14").Value = 0.82940000000000003

$ws.Range("G14").NumberFormat = $fmt164
$ws.Range("G14").Interior.Color = $colYellow
$ws.Range("G14").Value = 0.89139999999999997

$ws.Range("H14").NumberFormat = $fmt164
$ws.Range("H14").Interior.Color = $colRed
$ws.Range("H14").Value = 0.90090000000000003

# --- Row 15 ---
$ws.Range("B15").Value = "upper bound of CI for optimal F1 measure"
$ws.Range("C15").Value = "optimal F1 measure + critical t value * standard error of optimal F1 measure"

$ws.Range("D15").NumberFormat = $fmt164
$ws.Range("D15").Interior.Color = $colOrange
$ws.Range("D15").Value = 0.98809999999999998

$ws.Range("E15").NumberFormat = $fmt164
$ws.Range("E15").Interior.Color = $colGreen
$ws.Range("E15").Value = 0.98299999999999998

$ws.Range("F15").NumberFormat = $fmt164
$ws.Range("F15").Interior.Color = $colBlue
$ws.Range("F15").Value = 0.95299999999999996

$ws.Range("G15").NumberFormat = $fmt164
$ws.Range("G15").Interior.Color = $colYellow
$ws.Range("G15").Value = 0.98650000000000004

$ws.Range("H15").NumberFormat = $fmt164
$ws.Range("H15").Interior.Color = $colRed
$ws.Range("H15").Value = 0.99070000000000003

# Adjust column widths to fit the new, wider content (bestFit-style autosize)
$ws.Columns("B").ColumnWidth = 38.72135416666667
$ws.Columns("C").ColumnWidth = 71.05338541666667
$ws.Columns("D").ColumnWidth = 17.385416666666668
$ws.Columns("E").ColumnWidth = 22.608072916666668
$ws.Columns("F:H").ColumnWidth = 6.830729166666667

# Update the active selection to mirror the post-edit cursor position
$ws.Range("B18").Select() | Out-Null
